# "remove column from alcohol data"
# The alcohol measurement sheet (Sheet1) has a stray duplicate column: the
# old column M is dropped entirely and the data that used to live in column
# N shifts left to become the new column M.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(13).Delete() | Out-Null

# Leave the selection on the cell that now occupies the old spot (matches
# the author having clicked there after trimming the column).
$ws.Range("M1").Select() | Out-Null
